$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.06"
$ws.Range("E2").Value = "'0.91%"
$ws.Range("D3").Value = "'41.11"
$ws.Range("E3").Value = "'2.82%"
$ws.Range("D4").Value = "'5.732"
$ws.Range("E4").Value = "'0.56%"
$ws.Range("D5").Value = "'0.08100"
$ws.Range("E5").Value = "'0.22%"
$ws.Range("D6").Value = "'8.667"
$ws.Range("E6").Value = "'-0.18%"
$ws.Range("D7").Value = "'4.486"
$ws.Range("E7").Value = "'-1.82%"
$ws.Range("D8").Value = "'1.960"
$ws.Range("E8").Value = "'0.48%"
$ws.Range("D10").Value = "'0.9292"
$ws.Range("E10").Value = "'-1.50%"
$ws.Range("D11").Value = "'0.1263"
$ws.Range("E11").Value = "'-1.45%"
$ws.Range("D12").Value = "'0.1959"
$ws.Range("E12").Value = "'-1.22%"
$ws.Range("D13").Value = "'8.769"
$ws.Range("E13").Value = "'15.10%"
$ws.Range("D14").Value = "'0.09189"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("D15").Value = "'0.03731"
$ws.Range("E15").Value = "'9.10%"
$ws.Range("D16").Value = "'0.1050"
$ws.Range("E16").Value = "'9.33%"
$ws.Range("D17").Value = "'0.001293"
$ws.Range("E17").Value = "'-2.18%"
$ws.Range("D18").Value = "'0.006315"
$ws.Range("E18").Value = "'4.14%"
$ws.Range("D19").Value = "'3.367"
$ws.Range("E19").Value = "'-0.22%"
$ws.Range("D20").Value = "'0.3493"
$ws.Range("E20").Value = "'-0.37%"
$ws.Range("D21").Value = "'0.1365"
$ws.Range("E21").Value = "'-3.44%"
$ws.Range("D22").Value = "'0.2600"
$ws.Range("E22").Value = "'3.55%"
$ws.Range("D23").Value = "'0.04409"
$ws.Range("E23").Value = "'-0.40%"
$ws.Range("D24").Value = "'0.001251"
$ws.Range("E24").Value = "'-0.23%"
$ws.Range("D25").Value = "'0.004451"
$ws.Range("E25").Value = "'3.11%"
$ws.Range("D26").Value = "'0.0001238"
$ws.Range("E26").Value = "'3.82%"
$ws.Range("D39").Value = "'0.02743"
$ws.Range("E39").Value = "'8.50%"
$ws.Range("D40").Value = "'0.05542"
$ws.Range("E40").Value = "'6.58%"
$ws.Range("D41").Value = "'0.007536"
$ws.Range("E41").Value = "'3.22%"
$ws.Range("D42").Value = "'0.009812"
$ws.Range("E42").Value = "'8.48%"
$ws.Range("D43").Value = "'0.1425"
$ws.Range("E43").Value = "'-0.02%"
$ws.Range("D44").Value = "'0.002103"
$ws.Range("E44").Value = "'-4.13%"
$ws.Range("D45").Value = "'0.01181"
$ws.Range("E45").Value = "'17.81%"
$ws.Range("D46").Value = "'0.00006760"
$ws.Range("E46").Value = "'0.94%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.34%"
$ws.Range("D48").Value = "'0.003066"
$ws.Range("E48").Value = "'6.56%"
$ws.Range("D49").Value = "'0.002270"
$ws.Range("E49").Value = "'25.93%"
$ws.Range("D50").Value = "'0.00002096"
$ws.Range("E50").Value = "'-0.34%"
$ws.Range("D51").Value = "'0.0001996"
$ws.Range("E51").Value = "'-0.34%"
